# "Last fixes from me"
# Add the problem-size labels (10^6 / 10^7 / 10^8 / 10^9) next to each
# timing block on both sheets, and leave the workbook's active sheet /
# selection the way the author left it (Skalowalny active, Podstawowy's
# old H33 selection replaced by E25).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Podstawowy")
$ws2 = $wb.Worksheets.Item("Skalowalny")

# Podstawowy: one label per block of measurements (rows 2, 10, 18, 26)
$ws1.Range("D2").Value  = "10^6"
$ws1.Range("D10").Value = "10^7"
$ws1.Range("D18").Value = "10^8"
$ws1.Range("D26").Value = "10^9"

# Skalowalny: only the first block got the same treatment
$ws2.Range("D2").Value = "10^6"

# Final selection / active sheet, matching the saved workbook state
[void]$ws1.Range("E25").Select()
[void]$ws2.Activate()
[void]$ws2.Range("B1").Select()
